$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1332
$ws1.Range("F8").Value = 11543
$ws1.Range("F9").Value = 4356
$ws1.Range("F12").Value = 21
$ws1.Range("F14").Value = 2536
$ws1.Range("F16").Value = 134
$ws1.Range("F18").Value = 3614
$ws1.Range("F19").Value = 180
$ws1.Range("F21").Value = 11308
$ws1.Range("F22").Value = 11223
$ws1.Range("F24").Value = 42

# Sheet "全部类型" (sheet4): update "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1332
$ws4.Range("F8").Value = 11543
$ws4.Range("F9").Value = 4356
$ws4.Range("F12").Value = 21
$ws4.Range("F14").Value = 2536
$ws4.Range("F17").Value = 134
$ws4.Range("F19").Value = 3614
$ws4.Range("F20").Value = 180
$ws4.Range("F22").Value = 11308
$ws4.Range("F23").Value = 11223
$ws4.Range("F25").Value = 42
